$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Id (A), Ost (Q) and Nord (R) values between row 6 and row 7
$a6 = $ws.Range("A6").Value2
$a7 = $ws.Range("A7").Value2
$q6 = $ws.Range("Q6").Value2
$q7 = $ws.Range("Q7").Value2
$r6 = $ws.Range("R6").Value2
$r7 = $ws.Range("R7").Value2

$ws.Range("A6").Value2 = $a7
$ws.Range("A7").Value2 = $a6

$ws.Range("Q6").Value2 = $q7
$ws.Range("Q7").Value2 = $q6

$ws.Range("R6").Value2 = $r7
$ws.Range("R7").Value2 = $r6
